# resolved locator issues in addcandidate test
# Update the "Candidate" sheet: fix the contact_no test value and add a
# new "country" column (with value "USA") used by the updated locator.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Candidate")

# contact_no value fixed to a shorter/simpler test value
$ws.Range("D2").Value = 123

# new "country" column with header + value
$ws.Range("J1").Value = "country"
$ws.Range("J2").Value = "USA"

# match the author's recorded column widths (closest attainable via the
# Excel column-width property, which is quantized to 1/6-character steps)
$ws.Columns.Item(4).ColumnWidth = 10
$ws.Columns.Item(5).ColumnWidth = 9

# leave the selection on the newly added cell, as recorded in the workbook
$ws.Range("J2").Select()
